$d = $word.ActiveDocument

$ids = @("p110r_3", "p110v_1", "p110v_2", "p110v_3")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $rng = $d.Content
    $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)
}
